# Apply edits to "before.xlsx" per commit:
# "new sample report over 50 campaign members (14 campaign descriptions enhanced)"

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Campaign Data")
$ws2 = $wb.Worksheets.Item("Processing Summary")

# The COM engine stores ColumnWidth with a constant conversion offset relative to the
# OOXML "width" attribute (stored = ColumnWidth + 0.8333333333333333). Subtract the offset
# so the persisted <col width="..."> values match the target widths exactly.
$wOffset = 0.8333333333333333

# --- Resize columns B, C, F, L, M on "Campaign Data" ---
$ws1.Columns.Item(2).ColumnWidth = 38 - $wOffset
$ws1.Columns.Item(3).ColumnWidth = 16 - $wOffset
$ws1.Columns.Item(6).ColumnWidth = 19 - $wOffset
$ws1.Columns.Item(12).ColumnWidth = 24 - $wOffset
$ws1.Columns.Item(13).ColumnWidth = 36 - $wOffset

# --- Refresh the AI-generated sales descriptions for the existing 5 rows ---
$ws1.Range("W2").Value = 'Prospects were referred to Saasquatch via a trusted source, likely indicating strong interest or need. Their engagement suggests they are in the consideration stage of the buyer''s journey.'
$ws1.Range("W3").Value = 'Prospect self-submitted, likely urgently needing Invoca''s digital solution. Directly found us. Engaged with generic email, showing high initiative. In early to mid-buyer''s journey stage.'
$ws1.Range("W4").Value = 'Prospects actively searching for communication solutions with high intent found via RingCentral brand search, indicating IT decision makers in small businesses scaling up, likely in the consideration stage of the buyer''s journey.'
$ws1.Range("W5").Value = 'Prospect actively searched for communication solutions, focusing on RingCentral, showing high intent as a small business IT decision maker exploring UCaaS. Likely in the consideration stage of their buyer''s journey.'
$ws1.Range("W6").Value = 'Prospects actively searching for communication solutions, likely comparing options for a business phone system. Indicates mid-funnel buyer considering IT decisions for unified communications.'

# --- Append 9 new campaign rows (rows 7-15) ---
# Row 7
$ws1.Range("A7").Value = '7012H000001hVx6QAE'
$ws1.Range("B7").Value = 'Google_US_Search_Phone_Systems_Exact'
$ws1.Range("C7").Value = 'Paid Search'
$ws1.Range("D7").Value = 'Advertisement'
$ws1.Range("E7").Value = 'Paid Search SEM Soiurce - Everygreen Campaign Google_US_Search_Phone_Systems_Exact'
$ws1.Range("F7").Value = 'Search Engines'
$ws1.Range("G7").Value = 'SEM - Non-Brand'
$ws1.Range("H7").Value = 'MVP/PBX to Cloud'
$ws1.Range("I7").Value = 'RingEX'
$ws1.Range("J7").Value = 'Hello to Growing Your Business'
$ws1.Range("M7").Value = 'Google'
$ws1.Range("O7").Value = 'Easily set up and grow my business'
$ws1.Range("Q7").Value = 'US'
$ws1.Range("R7").Value = $false
$ws1.Range("T7").Value = 'Paid Search SEM Soiurce - Everygreen Campaign Google_US_Search_Phone_Systems_Exact'
$ws1.Range("U7").Value = 1
$ws1.Range("V7").Value = 'Based on the following campaign information, create a concise description (max 255 characters) that helps a salesperson understand:
            1. What the prospect was doing when they engaged with this campaign
            2. Why they likely engaged (their intent/interest)
            3. What this tells us about their buyer''s journey stage
            Focus on the prospect''s perspective and intent, not marketing terminology.
            IMPORTANT: If the campaign details mention any URLs or websites, preserve the domain name in your description.
            Campaign Information:
            Campaign: Google_US_Search_Phone_Systems_Exact
Engagement method: Clicked on search ads - actively searching for communication solutions with high intent
Cross channel marketing integration indicator: Migration from on-premise focus - cost savings pitch
Product interest: UCaaS/business phone system buyer - likely IT decision maker for unified communications
Secondary channel: SEM/SEO driven - paid or organic search
Specific engagement context: Searched generic terms like ''business phone system'' - comparing solutions
Target customer profile campaign identifier: Small business scaling - fast setup, flexible usage, and business growth
Campaign format: Ad campaign response - varying intent based on ad type
Lead source context: Google Ads - paid search intent
Value proposition focus: SMB messaging - simplicity and growth focus
Campaign description: Paid Search SEM Soiurce - Everygreen Campaign Google_US_Search_Phone_Systems_Exact
Campaign title: Google_US_Search_Phone_Systems_Exact
Target geographic market for campaign: US
Attribution tracking: Can clearly track that a lead came from this specific campaign (clear cause + effect)
Concise sales focused campaign summary: Paid Search SEM Soiurce - Everygreen Campaign Google_US_Search_Phone_Systems_Exact
            Description (max 255 characters):'
$ws1.Range("W7").Value = 'Prospects actively searching for business phone systems on Google US likely seek cost-effective UCaaS solutions. They are comparing options, indicating mid to late buyer''s journey stage as small businesses scaling up.'
$ws1.Rows.Item(7).RowHeight = 15

# Row 8
$ws1.Range("A8").Value = '7012H000001l35sQAA'
$ws1.Range("B8").Value = 'Affiliates_TheTop10sites'
$ws1.Range("C8").Value = 'Affiliates'
$ws1.Range("E8").Value = 'Prospect visited RingCentral Office Landing Page from a Comparison/Review Website - TheTop10Sites.com (https://www.thetop10sites.com/business-voip/) and submitted their info through the web form or called in.'
$ws1.Range("F8").Value = 'Direct Affiliates'
$ws1.Range("G8").Value = 'Affiliates - CPL'
$ws1.Range("H8").Value = 'MVP/PBX to Cloud'
$ws1.Range("I8").Value = 'RingEX'
$ws1.Range("K8").Value = 'Greenfield'
$ws1.Range("L8").Value = 'RingEX SMB Acquisition'
$ws1.Range("M8").Value = 'Better Impression (US)'
$ws1.Range("Q8").Value = 'US'
$ws1.Range("R8").Value = $false
$ws1.Range("T8").Value = 'Prospect visited RingCentral Office Landing Page from a Comparison/Review Website - TheTop10Sites.com (https://www.thetop10sites.com/business-voip/) and submitted their info through the web form or called in.'
$ws1.Range("U8").Value = 1
$ws1.Range("V8").Value = 'Based on the following campaign information, create a concise description (max 255 characters) that helps a salesperson understand:
            1. What the prospect was doing when they engaged with this campaign
            2. Why they likely engaged (their intent/interest)
            3. What this tells us about their buyer''s journey stage
            Focus on the prospect''s perspective and intent, not marketing terminology.
            IMPORTANT: If the campaign details mention any URLs or websites, preserve the domain name in your description.
            Campaign Information:
            Campaign: Affiliates_TheTop10sites
Engagement method: Referred through affiliate partner - has some context about RingCentral
Cross channel marketing integration indicator: Migration from on-premise focus - cost savings pitch
Product interest: UCaaS/business phone system buyer - likely IT decision maker for unified communications
Secondary channel: Referred by direct affiliate - warm intro with initial context
Specific engagement context: Affiliate referral - cost-per-lead, mid-intent form fill
Target customer profile program classification: Ongoing evergreen campaign - continous lead flow
Target customer profile and strategy: Targeting small business (1-499 employees) - faster sales cycle, price sensitive
Lead source context: Better Impression (US)
Company size segment: 20-99 employees - growing company, scalability important
Buyer journey stage: High intent - actively evaluating solutions (demo, trial, pricing interest)
Campaign description: Prospect visited RingCentral Office Landing Page from a Comparison/Review Website - TheTop10Sites.com (https://www.thetop10sites.com/business-voip/) and submitted their info through the web form or called in.
Campaign title: Affiliates_TheTop10sites
Target geographic market for campaign: US
Attribution tracking: Can clearly track that a lead came from this specific campaign (clear cause + effect)
Concise sales focused campaign summary: Prospect visited RingCentral Office Landing Page from a Comparison/Review Website - TheTop10Sites.com (https://www.thetop10sites.com/business-voip/) and submitted their info through the web form or called in.
            Description (max 255 characters):'
$ws1.Range("W8").Value = 'Prospect researching business phone systems on TheTop10Sites.com followed a cost-saving lead to RingCentral. High-intent evaluation stage, likely IT decision maker.'
$ws1.Rows.Item(8).RowHeight = 15

# Row 9
$ws1.Range("A9").Value = '7012H000001l3BBQAY'
$ws1.Range("B9").Value = 'Affiliates_BusinessBPS'
$ws1.Range("C9").Value = 'Affiliates'
$ws1.Range("E9").Value = 'Prospect visited RingCentral Office Landing Page from a Content/Review Website - Business.com (https://www.business.com/categories/business-phone-systems/) and submitted their info through the web form or called in.'
$ws1.Range("F9").Value = 'Direct Affiliates'
$ws1.Range("G9").Value = 'Affiliates - CPC'
$ws1.Range("H9").Value = 'MVP/PBX to Cloud'
$ws1.Range("I9").Value = 'RingEX'
$ws1.Range("K9").Value = 'Greenfield'
$ws1.Range("L9").Value = 'RingEX SMB Acquisition'
$ws1.Range("M9").Value = 'Purch (Business News Daily - VOIP)'
$ws1.Range("Q9").Value = 'US'
$ws1.Range("R9").Value = $false
$ws1.Range("T9").Value = 'Prospect visited RingCentral Office Landing Page from a Content/Review Website - Business.com (https://www.business.com/categories/business-phone-systems/) and submitted their info through the web form or called in.'
$ws1.Range("U9").Value = 1
$ws1.Range("V9").Value = 'Based on the following campaign information, create a concise description (max 255 characters) that helps a salesperson understand:
            1. What the prospect was doing when they engaged with this campaign
            2. Why they likely engaged (their intent/interest)
            3. What this tells us about their buyer''s journey stage
            Focus on the prospect''s perspective and intent, not marketing terminology.
            IMPORTANT: If the campaign details mention any URLs or websites, preserve the domain name in your description.
            Campaign Information:
            Campaign: Affiliates_BusinessBPS
Engagement method: Referred through affiliate partner - has some context about RingCentral
Cross channel marketing integration indicator: Migration from on-premise focus - cost savings pitch
Product interest: UCaaS/business phone system buyer - likely IT decision maker for unified communications
Secondary channel: Referred by direct affiliate - warm intro with initial context
Specific engagement context: Clicked cost-per-click affiliate link - low-friction awareness
Target customer profile program classification: Ongoing evergreen campaign - continous lead flow
Target customer profile and strategy: Targeting small business (1-499 employees) - faster sales cycle, price sensitive
Lead source context: Purch (Business News Daily - VOIP)
Company size segment: 20-99 employees - growing company, scalability important
Buyer journey stage: High intent - actively evaluating solutions (demo, trial, pricing interest)
Campaign description: Prospect visited RingCentral Office Landing Page from a Content/Review Website - Business.com (https://www.business.com/categories/business-phone-systems/) and submitted their info through the web form or called in.
Campaign title: Affiliates_BusinessBPS
Target geographic market for campaign: US
Attribution tracking: Can clearly track that a lead came from this specific campaign (clear cause + effect)
Concise sales focused campaign summary: Prospect visited RingCentral Office Landing Page from a Content/Review Website - Business.com (https://www.business.com/categories/business-phone-systems/) and submitted their info through the web form or called in.
            Description (max 255 characters):'
$ws1.Range("W9").Value = 'Prospect actively seeking unified communications solution for cost savings, likely an IT decision maker. High-intent evaluation stage, engaging through low-friction affiliate link from Business.com.'
$ws1.Rows.Item(9).RowHeight = 15

# Row 10
$ws1.Range("A10").Value = '70134000001CjkkAAC'
$ws1.Range("B10").Value = 'Bing_US_Search_Brand_Exact'
$ws1.Range("C10").Value = 'Paid Search'
$ws1.Range("E10").Value = 'Paid Search'
$ws1.Range("F10").Value = 'Brand Search'
$ws1.Range("G10").Value = 'SEM - Brand'
$ws1.Range("I10").Value = 'RingEX'
$ws1.Range("M10").Value = 'Bing'
$ws1.Range("Q10").Value = 'US'
$ws1.Range("R10").Value = $false
$ws1.Range("T10").Value = 'Paid Search'
$ws1.Range("U10").Value = 1
$ws1.Range("V10").Value = 'Based on the following campaign information, create a concise description (max 255 characters) that helps a salesperson understand:
            1. What the prospect was doing when they engaged with this campaign
            2. Why they likely engaged (their intent/interest)
            3. What this tells us about their buyer''s journey stage
            Focus on the prospect''s perspective and intent, not marketing terminology.
            IMPORTANT: If the campaign details mention any URLs or websites, preserve the domain name in your description.
            Campaign Information:
            Campaign: Bing_US_Search_Brand_Exact
Engagement method: Clicked on search ads - actively searching for communication solutions with high intent
Product interest: UCaaS/business phone system buyer - likely IT decision maker for unified communications
Secondary channel: Found via RingCentral brand keyword search - direct high intent query
Specific engagement context: Searched ''RingCentral'' or product names - brand aware, high intent
Lead source context: Bing
Campaign description: Paid Search
Campaign title: Bing_US_Search_Brand_Exact
Target geographic market for campaign: US
Attribution tracking: Can clearly track that a lead came from this specific campaign (clear cause + effect)
Concise sales focused campaign summary: Paid Search
            Description (max 255 characters):'
$ws1.Range("W10").Value = 'Prospects actively searching for communication solutions clicked on Bing search ads for ''RingCentral'' or related terms, showing high intent. Likely IT decision makers in the buying stage evaluating UCaaS options.'
$ws1.Rows.Item(10).RowHeight = 15

# Row 11
$ws1.Range("A11").Value = '70134000001XyCZAA0'
$ws1.Range("B11").Value = 'RCO/ACO_Price_Parity'
$ws1.Range("C11").Value = 'VAR MDF'
$ws1.Range("F11").Value = 'Events'
$ws1.Range("I11").Value = 'RingEX'
$ws1.Range("M11").Value = 'VAR Marketing'
$ws1.Range("Q11").Value = 'US'
$ws1.Range("R11").Value = $false
$ws1.Range("U11").Value = 1
$ws1.Range("V11").Value = 'Based on the following campaign information, create a concise description (max 255 characters) that helps a salesperson understand:
            1. What the prospect was doing when they engaged with this campaign
            2. Why they likely engaged (their intent/interest)
            3. What this tells us about their buyer''s journey stage
            Focus on the prospect''s perspective and intent, not marketing terminology.
            IMPORTANT: If the campaign details mention any URLs or websites, preserve the domain name in your description.
            Campaign Information:
            Campaign: RCO/ACO_Price_Parity
Engagement method: Campaign funded through reseller marketing - likely co-branded outreach via trusted tech advisor
Product interest: UCaaS/business phone system buyer - likely IT decision maker for unified communications
Secondary channel: Event marketing - webinars, conferences, tradeshows
Lead source context: VAR Marketing
Campaign title: RCO/ACO_Price_Parity
Target geographic market for campaign: US
Attribution tracking: Can clearly track that a lead came from this specific campaign (clear cause + effect)
            Description (max 255 characters):'
$ws1.Range("W11").Value = 'Prospects engaged with the RCO/ACO_Price_Parity campaign, seeking price parity for UCaaS solutions. Likely IT decision makers exploring unified communications options, indicating a mid-stage buyer''s journey.'
$ws1.Rows.Item(11).RowHeight = 15

# Row 12
$ws1.Range("A12").Value = '70180000000OwaeAAC'
$ws1.Range("B12").Value = 'SEO_GoogleRC'
$ws1.Range("C12").Value = 'Organic Search'
$ws1.Range("F12").Value = 'Search Engines'
$ws1.Range("I12").Value = 'RingEX'
$ws1.Range("Q12").Value = 'US'
$ws1.Range("R12").Value = $false
$ws1.Range("U12").Value = 5
$ws1.Range("V12").Value = 'Based on the following campaign information, create a concise description (max 255 characters) that helps a salesperson understand:
            1. What the prospect was doing when they engaged with this campaign
            2. Why they likely engaged (their intent/interest)
            3. What this tells us about their buyer''s journey stage
            Focus on the prospect''s perspective and intent, not marketing terminology.
            IMPORTANT: If the campaign details mention any URLs or websites, preserve the domain name in your description.
            Campaign Information:
            Campaign: SEO_GoogleRC
Engagement method: Found us through organic search - self-directed research, comparing options
Product interest: UCaaS/business phone system buyer - likely IT decision maker for unified communications
Secondary channel: SEM/SEO driven - paid or organic search
Campaign title: SEO_GoogleRC
Target geographic market for campaign: US
Attribution tracking: Can clearly track that a lead came from this specific campaign (clear cause + effect)
            Description (max 255 characters):'
$ws1.Range("W12").Value = 'Prospects conducting self-research on Google for UCaaS options in the US found us through SEO_GoogleRC. Likely IT decision makers at the comparison stage of their buyer''s journey.'
$ws1.Rows.Item(12).RowHeight = 15

# Row 13
$ws1.Range("A13").Value = '701800000019F0iAAE'
$ws1.Range("B13").Value = 'Sales Generated'
$ws1.Range("C13").Value = 'Walk-On'
$ws1.Range("F13").Value = 'Sales Generated'
$ws1.Range("I13").Value = 'General'
$ws1.Range("Q13").Value = 'US'
$ws1.Range("R13").Value = $false
$ws1.Range("U13").Value = 21
$ws1.Range("V13").Value = 'Based on the following campaign information, create a concise description (max 255 characters) that helps a salesperson understand:
            1. What the prospect was doing when they engaged with this campaign
            2. Why they likely engaged (their intent/interest)
            3. What this tells us about their buyer''s journey stage
            Focus on the prospect''s perspective and intent, not marketing terminology.
            IMPORTANT: If the campaign details mention any URLs or websites, preserve the domain name in your description.
            Campaign Information:
            Campaign: Sales Generated
Engagement method: Self-submitted or inbound lead without campaign - high initiative, potentially urgent need
Secondary channel: Entered by sales team from internal source - low initial engagement
Campaign title: Sales Generated
Target geographic market for campaign: US
Attribution tracking: Can clearly track that a lead came from this specific campaign (clear cause + effect)
            Description (max 255 characters):'
$ws1.Range("W13").Value = 'Prospects likely sought help urgently when engaging with "Sales Generated" campaign from the US. Their high initiative suggests immediate needs, placing them in the decision-making stage of the buyer''s journey.'
$ws1.Rows.Item(13).RowHeight = 15

# Row 14
$ws1.Range("A14").Value = '701Hr000000t6JrIAI'
$ws1.Range("B14").Value = 'HealthcareLPformfills'
$ws1.Range("C14").Value = 'Other'
$ws1.Range("E14").Value = 'Healthcare integrated campaign'
$ws1.Range("F14").Value = 'website'
$ws1.Range("I14").Value = 'General'
$ws1.Range("M14").Value = 'Healthcare'
$ws1.Range("N14").Value = 'Healthcare'
$ws1.Range("Q14").Value = 'US'
$ws1.Range("R14").Value = $false
$ws1.Range("T14").Value = 'Healthcare integrated campaign'
$ws1.Range("U14").Value = 1
$ws1.Range("V14").Value = 'Based on the following campaign information, create a concise description (max 255 characters) that helps a salesperson understand:
            1. What the prospect was doing when they engaged with this campaign
            2. Why they likely engaged (their intent/interest)
            3. What this tells us about their buyer''s journey stage
            Focus on the prospect''s perspective and intent, not marketing terminology.
            IMPORTANT: If the campaign details mention any URLs or websites, preserve the domain name in your description.
            Campaign Information:
            Campaign: HealthcareLPformfills
Engagment channel not categorized - review needed to determine buyer intent
Secondary channel: Direct website visit or form fill - proactive buyer behavior
Lead source context: Healthcare
Industry context: Healthcare industry - HIPAA compilance needs
Campaign description: Healthcare integrated campaign
Campaign title: HealthcareLPformfills
Target geographic market for campaign: US
Attribution tracking: Can clearly track that a lead came from this specific campaign (clear cause + effect)
Concise sales focused campaign summary: Healthcare integrated campaign
            Description (max 255 characters):'
$ws1.Range("W14").Value = 'Prospects actively sought healthcare solutions via direct website visits or form fills, showing proactive interest in HIPAA compliance. Indicates mid to late buyer''s journey stages.'
$ws1.Rows.Item(14).RowHeight = 15

# Row 15
$ws1.Range("A15").Value = '701Hr000002I3SJIA0'
$ws1.Range("B15").Value = 'True_Walk_On_2024'
$ws1.Range("C15").Value = 'Walk-On'
$ws1.Range("F15").Value = 'Walk-On'
$ws1.Range("I15").Value = 'General'
$ws1.Range("Q15").Value = 'US'
$ws1.Range("R15").Value = $false
$ws1.Range("U15").Value = 5
$ws1.Range("V15").Value = 'Based on the following campaign information, create a concise description (max 255 characters) that helps a salesperson understand:
            1. What the prospect was doing when they engaged with this campaign
            2. Why they likely engaged (their intent/interest)
            3. What this tells us about their buyer''s journey stage
            Focus on the prospect''s perspective and intent, not marketing terminology.
            IMPORTANT: If the campaign details mention any URLs or websites, preserve the domain name in your description.
            Campaign Information:
            Campaign: True_Walk_On_2024
Engagement method: Self-submitted or inbound lead without campaign - high initiative, potentially urgent need
Secondary channel: Inbound lead without source - prospect found us directly
Campaign title: True_Walk_On_2024
Target geographic market for campaign: US
Attribution tracking: Can clearly track that a lead came from this specific campaign (clear cause + effect)
            Description (max 255 characters):'
$ws1.Range("W15").Value = 'The prospect proactively sought information about "True_Walk_On_2024," indicating a high initiative and potentially urgent need. They likely engaged due to a direct interest in the campaign topic, showing an advanced stage in their buyer''s journey.'
$ws1.Rows.Item(15).RowHeight = 15

# --- Update "Processing Summary" metrics to reflect the new totals ---
$ws2.Range("B3").Value = 14
$ws2.Range("B4").Value = 14
$ws2.Range("B7").Value = '201.8 chars'
$ws2.Range("B8").Value = 50
$ws2.Range("B9").Value = 0.35
$ws2.Range("B10").Value = 7
$ws2.Range("B11").Value = 1
$ws2.Range("B13").Value = 13
$ws2.Range("B15").Value = 14
$ws2.Range("B16").Value = 9
$ws2.Range("B17").Value = '2025-07-14 12:31:06'
